# Apply the weekly report refresh edit:
#  - bump the "Report Generated On" timestamp
#  - populate the billed-amount / total figures now that the line item has pricing
#  - clear the now-unused Scope ID value in G10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report generation timestamp (D5)
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:00 AM"

# Total Billed Amount (C8) and the corresponding line-item / total pricing (H16, H17)
$ws.Range("C8").Value = 478.55
$ws.Range("H16").Value = 478.55
$ws.Range("H17").Value = 478.55

# Scope ID #: value cleared out
$ws.Range("G10").Value = ""
